$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(1.02, 1.068669265048171, 1.070453912345782, 1.070370114642192, 1.081441882673263, 1.058013189309711, 1.073607896811492, 1.07315370785997, 1.073070133654963, 1.084112696005698, 1.075132543169511),
  @(1.02, 1.069899380881899, 1.071438181874338, 1.07150401731068, 1.082606192036813, 1.058443955918806, 1.074493543146562, 1.073953926814216, 1.074019599382579, 1.085094613439919, 1.076019447223962),
  @(1.02, 1.070694910225125, 1.072074660466521, 1.072237900776096, 1.083359658404341, 1.058721165755449, 1.075065622461784, 1.074470693502606, 1.074633550351382, 1.085729474036999, 1.076592338957431),
  @(1.02, 1.071029248169691, 1.072342139121407, 1.07254646792852, 1.083676436253126, 1.058837340457921, 1.075305888250761, 1.074687697151383, 1.074891556687947, 1.085996249633062, 1.076832945951387),
  @(1.02, 1.071085379030624, 1.072387044312013, 1.072598280241405, 1.083729625894382, 1.058856825360458, 1.075346216114396, 1.074724118681557, 1.074934871335469, 1.086041035440857, 1.076873331085214),
  @(1.02, 1.070699378065424, 1.072078234908936, 1.072242023696795, 1.083363891122874, 1.058722719518055, 1.075068833831995, 1.074473594077489, 1.074636998230879, 1.085733039174734, 1.076595554888155),
  @(1.02, 1.069085080458847, 1.070786636326255, 1.070753287349999, 1.081835350649703, 1.058159085117703, 1.073907411600023, 1.073424359251934, 1.073391097041646, 1.084444645050827, 1.075432483303397),
  @(1.02, 1.066237028173789, 1.06850748705325, 1.068131210865493, 1.079142421593573, 1.057154175347892, 1.071853168756838, 1.071567544455901, 1.071192415985923, 1.082170400954015, 1.07337532319975),
  @(1.02, 1.06433586232277, 1.066985838034162, 1.066383932700929, 1.077347420745972, 1.056476310044812, 1.070478427967093, 1.070324262126368, 1.069724373503026, 1.080651515375335, 1.071998630120398),
  @(1.02, 1.063512019654623, 1.066326404514256, 1.065627506891078, 1.076570212893297, 1.056180895575066, 1.069881885210077, 1.069784607662352, 1.069088142861996, 1.079993156536063, 1.071401240204243),
  @(1.02, 1.063205910975523, 1.066081377426144, 1.065346558600178, 1.076281527356194, 1.05607087982551, 1.069660109853962, 1.069583958320107, 1.068851732988328, 1.079748510170806, 1.071179149901684),
  @(1.02, 1.063271576754486, 1.066133940400055, 1.06540682199666, 1.076343451229682, 1.056094491521138, 1.069707690176108, 1.069627007254515, 1.068902447602349, 1.079800992316116, 1.071226797793345),
  @(1.02, 1.063486718607741, 1.066306152223781, 1.065604283179961, 1.076546349984538, 1.056171807472459, 1.069863557134179, 1.069768025966223, 1.069068602901045, 1.079972936093702, 1.071382886100375),
  @(1.02, 1.063619261686714, 1.066412246418682, 1.065725948388371, 1.076671363087122, 1.056219406497008, 1.069959566305317, 1.069854886011247, 1.069170965289794, 1.080078862685877, 1.071479031615548),
  @(1.02, 1.064390524626758, 1.067029590755002, 1.066434137358874, 1.077399002143052, 1.056495875728497, 1.070517991608543, 1.070360049585732, 1.069766586133298, 1.080695194192193, 1.072038249946756),
  @(1.02, 1.064874148044303, 1.067416686336967, 1.066878406483642, 1.077855440132257, 1.056668789660496, 1.070867935464244, 1.070676574844573, 1.070140052703352, 1.08108162171082, 1.072388690762872),
  @(1.02, 1.065156177264611, 1.067642419733888, 1.06713755642656, 1.078121676613023, 1.056769464729564, 1.071071929296823, 1.070861072866339, 1.070357835649428, 1.081306953685148, 1.072592974290092),
  @(1.02, 1.065252331826896, 1.067719380061748, 1.067225922501852, 1.078212457093739, 1.056803761368862, 1.07114146520049, 1.070923960632382, 1.07043208493003, 1.081383775149302, 1.072662608942718),
  @(1.02, 1.064822266070927, 1.067375160129163, 1.066830739043878, 1.077806468264275, 1.056650256538007, 1.070830402485392, 1.070642627704937, 1.070099988850418, 1.081040168396591, 1.072351104482885),
  @(1.02, 1.063423367390427, 1.066255442468642, 1.065546135191052, 1.076486601261022, 1.056149047746057, 1.069817663537619, 1.069726504939162, 1.069019676679314, 1.07992230578226, 1.071336927329644),
  @(1.02, 1.062543262602969, 1.065550944074994, 1.064738579837795, 1.075656771088066, 1.055832264605157, 1.069179797865827, 1.069149357888624, 1.068339946169444, 1.07921886684688, 1.070698155815433),
  @(1.02, 1.063009877007665, 1.065924458733983, 1.065166668820774, 1.076096678240546, 1.056000354409064, 1.069518048985626, 1.069455423457128, 1.068700331653536, 1.079591830046601, 1.07103688729062),
  @(1.02, 1.064845709487257, 1.067393924198544, 1.066852277860971, 1.077828596531404, 1.056658631422845, 1.070847362401469, 1.07065796734486, 1.070118092147288, 1.081058899564801, 1.072368088483988),
  @(1.02, 1.066973741864891, 1.069097087013329, 1.068808938463069, 1.079838551597708, 1.057415361879966, 1.072385157746879, 1.072048522900789, 1.071761219787197, 1.082758821524309, 1.073908067675176)
)

$colLetters = @("B","C","D","E","F","I","J","K","L","M","N")
$startRow = 2

for ($r = 0; $r -lt $data.Count; $r++) {
    for ($c = 0; $c -lt $colLetters.Count; $c++) {
        $cellRef = "$($colLetters[$c])$($startRow + $r)"
        $ws.Range($cellRef).Value = $data[$r][$c]
    }
}